$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 62
$ws_ALC.Range("H62").Value = 11998.333
$ws_ALC.Range("I62").Value = 10000
$ws_ALC.Range("J62").Value = 12997.5
$ws_ALC.Range("K62").Value = 10000
$ws_ALC.Range("L62").Value = 12997.5
$ws_ALC.Range("M62").Value = -9376
$ws_ALC.Range("N62").Value = -14245.5

# ALC row 65
$ws_ALC.Range("H65").Value = 11998.333
$ws_ALC.Range("I65").Value = 10000
$ws_ALC.Range("J65").Value = 12997.5
$ws_ALC.Range("K65").Value = 50000
$ws_ALC.Range("L65").Value = 64987.5
$ws_ALC.Range("M65").Value = -46880
$ws_ALC.Range("N65").Value = -71227.5

# ALC row 76
$ws_ALC.Range("H76").Value = 5748
$ws_ALC.Range("I76").Value = 5897.6
$ws_ALC.Range("J76").Value = 5000
$ws_ALC.Range("K76").Value = 5897.6
$ws_ALC.Range("L76").Value = 5000
$ws_ALC.Range("M76").Value = -5582.6
$ws_ALC.Range("N76").Value = -5630

# ALC row 79
$ws_ALC.Range("H79").Value = 5748
$ws_ALC.Range("I79").Value = 5897.6
$ws_ALC.Range("J79").Value = 5000
$ws_ALC.Range("K79").Value = 5897.6
$ws_ALC.Range("L79").Value = 5000
$ws_ALC.Range("M79").Value = -4805.6
$ws_ALC.Range("N79").Value = -7184

# ALC row 86
$ws_ALC.Range("H86").Value = 0
$ws_ALC.Range("I86").Value = 0
$ws_ALC.Range("J86").Value = 0
$ws_ALC.Range("K86").Value = 0
$ws_ALC.Range("L86").Value = 0
$ws_ALC.Range("M86").ClearContents()
$ws_ALC.Range("N86").ClearContents()

# ALC row 87
$ws_ALC.Range("H87").Value = 0
$ws_ALC.Range("J87").Value = 0
$ws_ALC.Range("L87").Value = 0
$ws_ALC.Range("N87").ClearContents()

# ALC row 89
$ws_ALC.Range("H89").Value = 0
$ws_ALC.Range("I89").Value = 0
$ws_ALC.Range("J89").Value = 0
$ws_ALC.Range("K89").Value = 0
$ws_ALC.Range("L89").Value = 0
$ws_ALC.Range("M89").ClearContents()
$ws_ALC.Range("N89").ClearContents()

# ALC row 90
$ws_ALC.Range("H90").Value = 0
$ws_ALC.Range("J90").Value = 0
$ws_ALC.Range("L90").Value = 0
$ws_ALC.Range("N90").ClearContents()

# ALC row 98
$ws_ALC.Range("H98").Value = 7066.909
$ws_ALC.Range("I98").Value = 1882.3
$ws_ALC.Range("J98").Value = 11387.417
$ws_ALC.Range("K98").Value = 1882.3
$ws_ALC.Range("L98").Value = 11387.417
$ws_ALC.Range("M98").Value = -384.3
$ws_ALC.Range("N98").Value = -14383.417

# ALC row 122
$ws_ALC.Range("H122").Value = 7066.909
$ws_ALC.Range("I122").Value = 1882.3
$ws_ALC.Range("J122").Value = 11387.417
$ws_ALC.Range("K122").Value = 5646.9
$ws_ALC.Range("L122").Value = 34162.251
$ws_ALC.Range("M122").Value = -3196.9
$ws_ALC.Range("N122").Value = -39062.251

# ALC row 138
$ws_ALC.Range("H138").Value = 5249.5
$ws_ALC.Range("I138").Value = 1999.3334
$ws_ALC.Range("J138").Value = 15000
$ws_ALC.Range("K138").Value = 5998.0002
$ws_ALC.Range("L138").Value = 45000
$ws_ALC.Range("M138").Value = -858.0002000000004
$ws_ALC.Range("N138").Value = -55280

# ARM row 5
$ws_ARM.Range("H5").Value = 1
$ws_ARM.Range("J5").Value = 1
$ws_ARM.Range("L5").Value = 1
$ws_ARM.Range("N5").Value = -225

# ARM row 11
$ws_ARM.Range("H11").Value = 289569.56
$ws_ARM.Range("I11").Value = 669333.3
$ws_ARM.Range("J11").Value = 4746.75
$ws_ARM.Range("K11").Value = 669333.3
$ws_ARM.Range("L11").Value = 4746.75
$ws_ARM.Range("M11").Value = -669189.3
$ws_ARM.Range("N11").Value = -5034.75

# ARM row 74
$ws_ARM.Range("H74").Value = 8283.857
$ws_ARM.Range("J74").Value = 9253.5
$ws_ARM.Range("L74").Value = 9253.5
$ws_ARM.Range("N74").Value = -11001.5

# ARM row 77
$ws_ARM.Range("H77").Value = 8283.857
$ws_ARM.Range("J77").Value = 9253.5
$ws_ARM.Range("L77").Value = 46267.5
$ws_ARM.Range("N77").Value = -55003.5

# ARM row 111
$ws_ARM.Range("H111").Value = 52644
$ws_ARM.Range("J111").Value = 52644
$ws_ARM.Range("L111").Value = 52644
$ws_ARM.Range("N111").Value = -60824

# BSM row 4
$ws_BSM.Range("H4").Value = 1
$ws_BSM.Range("J4").Value = 1
$ws_BSM.Range("L4").Value = 1
$ws_BSM.Range("N4").Value = -231

# CRP row 13
$ws_CRP.Range("H13").Value = 554
$ws_CRP.Range("J13").Value = 554
$ws_CRP.Range("L13").Value = 554
$ws_CRP.Range("N13").Value = -832

# CRP row 58
$ws_CRP.Range("H58").Value = 8970.450999999999
$ws_CRP.Range("I58").Value = 832.55554
$ws_CRP.Range("K58").Value = 832.55554
$ws_CRP.Range("M58").Value = -629.55554

# CRP row 105
$ws_CRP.Range("H105").Value = 7669.6665
$ws_CRP.Range("I105").Value = 14998
$ws_CRP.Range("K105").Value = 14998
$ws_CRP.Range("M105").Value = -13251

# CRP row 134
$ws_CRP.Range("H134").Value = 7261.8887
$ws_CRP.Range("I134").Value = 4479.5713
$ws_CRP.Range("K134").Value = 13438.7139
$ws_CRP.Range("M134").Value = -10903.7139

# CRP row 136
$ws_CRP.Range("H136").Value = 8970.450999999999
$ws_CRP.Range("I136").Value = 832.55554
$ws_CRP.Range("K136").Value = 2497.66662
$ws_CRP.Range("M136").Value = 52.33338000000003

# CUL row 117
$ws_CUL.Range("H117").Value = 4209.6
$ws_CUL.Range("J117").Value = 4682.6665
$ws_CUL.Range("L117").Value = 14047.9995
$ws_CUL.Range("N117").Value = -20931.9995

# CUL row 129
$ws_CUL.Range("H129").Value = 1048.75
$ws_CUL.Range("I129").Value = 1065
$ws_CUL.Range("K129").Value = 3195
$ws_CUL.Range("M129").Value = 1805

# GSM row 4
$ws_GSM.Range("H4").Value = 3000
$ws_GSM.Range("J4").Value = 3000
$ws_GSM.Range("L4").Value = 3000
$ws_GSM.Range("N4").Value = -3224

# GSM row 70
$ws_GSM.Range("H70").Value = 7499.5
$ws_GSM.Range("I70").Value = 7499.5
$ws_GSM.Range("K70").Value = 7499.5
$ws_GSM.Range("M70").Value = -7229.5

# GSM row 73
$ws_GSM.Range("H73").Value = 7499.5
$ws_GSM.Range("I73").Value = 7499.5
$ws_GSM.Range("K73").Value = 7499.5
$ws_GSM.Range("M73").Value = -6563.5

# GSM row 80
$ws_GSM.Range("H80").Value = 2435
$ws_GSM.Range("I80").Value = 2652.5
$ws_GSM.Range("J80").Value = 2000
$ws_GSM.Range("K80").Value = 2652.5
$ws_GSM.Range("L80").Value = 2000
$ws_GSM.Range("M80").Value = -1654.5
$ws_GSM.Range("N80").Value = -3996

# GSM row 83
$ws_GSM.Range("H83").Value = 2435
$ws_GSM.Range("I83").Value = 2652.5
$ws_GSM.Range("J83").Value = 2000
$ws_GSM.Range("K83").Value = 13262.5
$ws_GSM.Range("L83").Value = 10000
$ws_GSM.Range("M83").Value = -8270.5
$ws_GSM.Range("N83").Value = -19984

# GSM row 118
$ws_GSM.Range("H118").Value = 26102.666
$ws_GSM.Range("J118").Value = 26102.666
$ws_GSM.Range("L118").Value = 26102.666
$ws_GSM.Range("N118").Value = -29416.666

# GSM row 132
$ws_GSM.Range("H132").Value = 6366.8667
$ws_GSM.Range("I132").Value = 3980.4
$ws_GSM.Range("K132").Value = 11941.2
$ws_GSM.Range("M132").Value = -9411.200000000001

# LTW row 5
$ws_LTW.Range("H5").Value = 6999
$ws_LTW.Range("I5").Value = 6999
$ws_LTW.Range("K5").Value = 6999
$ws_LTW.Range("M5").Value = -6886

# LTW row 16
$ws_LTW.Range("H16").Value = 14998
$ws_LTW.Range("I16").Value = 14998
$ws_LTW.Range("K16").Value = 14998
$ws_LTW.Range("M16").Value = -14828

# LTW row 40
$ws_LTW.Range("H40").Value = 8124.5
$ws_LTW.Range("I40").Value = 8124.5
$ws_LTW.Range("K40").Value = 8124.5
$ws_LTW.Range("M40").Value = -7988.5

# LTW row 43
$ws_LTW.Range("H43").Value = 0
$ws_LTW.Range("J43").Value = 0
$ws_LTW.Range("L43").Value = 0
$ws_LTW.Range("N43").ClearContents()

# LTW row 53
$ws_LTW.Range("H53").Value = 2000
$ws_LTW.Range("I53").Value = 2000
$ws_LTW.Range("K53").Value = 2000
$ws_LTW.Range("M53").Value = -1482

# LTW row 96
$ws_LTW.Range("H96").Value = 69999
$ws_LTW.Range("J96").Value = 69999
$ws_LTW.Range("L96").Value = 69999
$ws_LTW.Range("N96").Value = -75491

# LTW row 132
$ws_LTW.Range("H132").Value = 9610.625
$ws_LTW.Range("I132").Value = 7277
$ws_LTW.Range("K132").Value = 21831
$ws_LTW.Range("M132").Value = -19301

# LTW row 136
$ws_LTW.Range("H136").Value = 11672.842
$ws_LTW.Range("I136").Value = 6141.5713
$ws_LTW.Range("J136").Value = 14899.417
$ws_LTW.Range("K136").Value = 18424.7139
$ws_LTW.Range("L136").Value = 44698.251
$ws_LTW.Range("M136").Value = -15874.7139
$ws_LTW.Range("N136").Value = -49798.251

# WVR row 74
$ws_WVR.Range("H74").Value = 46733.6
$ws_WVR.Range("J74").Value = 47167
$ws_WVR.Range("L74").Value = 47167
$ws_WVR.Range("N74").Value = -49039

# WVR row 77
$ws_WVR.Range("H77").Value = 46733.6
$ws_WVR.Range("J77").Value = 47167
$ws_WVR.Range("L77").Value = 141501
$ws_WVR.Range("N77").Value = -150861
